$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,7,8,11,14,15)  # B,C,D,E,G,H,K,N,O column indices
$data = @(
  @(9.387601030044342, 6.137134154901638, 5.975149884727422, 16.49042890569262, 26.37899438348488, 13.59947067343654, 8.534746551802261, 17.61626620208107, 20.39876267414844),
  @(9.046875586347332, 5.974529650032757, 5.854373674938505, 15.55648342679545, 26.34056623700892, 13.63902358627923, 8.290465960405923, 17.6764094635205, 20.44519499119439),
  @(8.832643403211399, 5.871489701444482, 5.780693619828576, 14.95817113337863, 26.32653548527147, 13.66569628805197, 8.13815724101206, 17.7150116354277, 20.47862416110296),
  @(8.744219215591205, 5.828732981768979, 5.75083538619954, 14.70837106189489, 26.32322260302507, 13.67716512111848, 8.075603316382059, 17.7311646316885, 20.49348005485504),
  @(8.729473011169109, 5.821588084152005, 5.745888920126692, 14.66653952453227, 26.32281768034569, 13.67910569742263, 8.065189810899914, 17.7338723702928, 20.49602123173506),
  @(8.831455233979636, 5.870916123573141, 5.780290202552083, 14.95482607164587, 26.32648107230253, 13.66584853446636, 8.137315456586187, 17.71522776872179, 20.4788195248665),
  @(9.271242241606782, 6.081753978112168, 5.933435254004691, 16.17371166815139, 26.36375998841913, 13.61261263854785, 8.451055935566448, 17.63665699288009, 20.41374927134245),
  @(10.08828167396873, 6.468309155665016, 6.235512026630915, 18.44066520877717, 26.51263334769786, 13.52719060672823, 9.048931511905135, 17.49579869830934, 20.32534071682669),
  @(10.6544369669269, 6.734098729559121, 6.455897845087441, 20.07755669194932, 26.66785291661464, 13.47604238890276, 9.600766863175583, 17.40027826761627, 20.2844819668715),
  @(10.90346728911329, 6.850729474976704, 6.55529938777506, 20.78000746192793, 26.74828532983003, 13.4553041698097, 9.839937131234198, 17.35853442758631, 20.27116259530068),
  @(10.99646568254752, 6.894255205835051, 6.592776948309719, 21.03997782133113, 26.78013948943381, 13.44781558584727, 9.928768486639933, 17.3429714212514, 20.26687870789882),
  @(10.97649596553517, 6.884909981869767, 6.584713390519238, 20.98425632962147, 26.77321732082226, 13.44941216128589, 9.909714771982628, 17.3463123397026, 20.26776748998471),
  @(10.91114487513704, 6.85432331778761, 6.558386195238636, 20.80151609290564, 26.75087811292178, 13.45468076926255, 9.84728031786754, 17.35724915710647, 20.27079491594441),
  @(10.87094344692001, 6.835504093832909, 6.542237540588395, 20.68879778224404, 26.73737599230311, 13.45795543976603, 9.808810332780874, 17.36398007495749, 20.27274832421904),
  @(10.63798386095073, 6.72638838792023, 6.44938099826297, 20.03080393792604, 26.66279279067518, 13.47744864753702, 9.584895008184899, 17.40304060530785, 20.28545862925721),
  @(10.49282969190743, 6.658335385580188, 6.392168131441444, 19.61636631829482, 26.6195434394749, 13.49005551629083, 9.444467814059509, 17.42743976535724, 20.29460695536805),
  @(10.40854329206466, 6.618791716837282, 6.359182815451395, 19.37402069004431, 26.5955932289918, 13.49754468272275, 9.362582760433996, 17.44163443943886, 20.30036451013368),
  @(10.3798709687226, 6.605334762182577, 6.348002446486273, 19.29128341282431, 26.58764353437136, 13.50012124230221, 9.334667427266368, 17.44646818921292, 20.3023989729598),
  @(10.50836475243058, 6.665621471312296, 6.398266914565749, 19.66089476621191, 26.62405171207135, 13.48868885248745, 9.459532176786775, 17.42482578733873, 20.29358178072542),
  @(10.93037604634622, 6.863324910343989, 6.5661238868044, 20.8553547899133, 26.75740193220605, 13.45312335186449, 9.865666202170214, 17.35403012356668, 20.26988504759822),
  @(11.19855434678827, 6.988796574896819, 6.674856952235745, 21.60085650105215, 26.8526823480772, 13.43200435055242, 10.12096143296773, 17.30918550468537, 20.25882780811433),
  @(11.05614436491192, 6.922179600218602, 6.616925851153852, 21.20617336672747, 26.8010917457984, 13.44308125230065, 9.985642100297323, 17.33298999787495, 20.26432322794069),
  @(10.50134394730413, 6.662328735492, 6.395509940878433, 19.64077616600321, 26.62201067238685, 13.48930596955053, 9.452725162967814, 17.42600704562835, 20.29404371087663),
  @(9.872824925714733, 6.366813442146914, 6.153880896025936, 17.80024299104303, 26.46427196069235, 13.54826358349575, 8.886486882389672, 17.53249881378627, 20.34503948469115)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $c = $cols[$j]
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}

Write-Output "done"